# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> stock "Office Theme" palette (unused/orphaned - only
#                             wired from the notes master relationships)
#   ppt/theme/theme2.xml  -> "Integral" palette, the theme actually applied to the
#                             slide master / presentation (this is what renders)
# The authored edit swaps the two palettes: the live design ("Integral") becomes
# the stock "Office" colours. Font scheme and format scheme are identical between
# the two theme parts, so only the 12 theme colours need to move.
#
# Walk to the live theme's colour scheme through the Design -> SlideMaster -> Theme
# chain and overwrite each of the 12 slots (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) with the values the "Office Theme" palette used.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function Set-ThemeRGB {
    param(
        [int]$Index,
        [string]$HexColor
    )
    $r = [Convert]::ToInt32($HexColor.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($HexColor.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($HexColor.Substring(4, 2), 16)
    $comRgb = $r + ($g * 256) + ($b * 65536)
    $colorScheme.Colors($Index).RGB = $comRgb
}

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeRGB 1  "000000"
Set-ThemeRGB 2  "FFFFFF"
Set-ThemeRGB 3  "44546A"
Set-ThemeRGB 4  "E7E6E6"
Set-ThemeRGB 5  "5B9BD5"
Set-ThemeRGB 6  "ED7D31"
Set-ThemeRGB 7  "A5A5A5"
Set-ThemeRGB 8  "FFC000"
Set-ThemeRGB 9  "4472C4"
Set-ThemeRGB 10 "70AD47"
Set-ThemeRGB 11 "0563C1"
Set-ThemeRGB 12 "954F72"
